# Update "想去人数" (want-to-go count) values on the 展览 and 全部类型 sheets.
# These values were regenerated by the gh-pages build at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 180
$wsExhibit.Range("F5").Value = 26
$wsExhibit.Range("F11").Value = 1880
$wsExhibit.Range("F20").Value = 36
$wsExhibit.Range("F21").Value = 40
$wsExhibit.Range("F23").Value = 1014
$wsExhibit.Range("F27").Value = 254

# Sheet "全部类型" (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 180
$wsAll.Range("F5").Value = 26
$wsAll.Range("F12").Value = 1880
$wsAll.Range("F21").Value = 36
$wsAll.Range("F22").Value = 40
$wsAll.Range("F24").Value = 1014
$wsAll.Range("F28").Value = 254
